# "Laporan Buku Tamu" (guest book) update:
# Remove the 2024-09-05 / Sion / Ikebukuro guest-book entry (row 5) and the
# last three entries (old rows 11-13: Rofi/otto-iskandar, Albert/Jakarta,
# Rafa/Warungkondang), leaving the table with 8 guest entries instead of 12.
# Remaining rows shift up and the "No" column is renumbered sequentially.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 5 (2024-09-05, Sion, Ikebukuro, ...). Rows below shift up.
$ws.Rows.Item(5).Delete()

# After the shift, the former rows 11-13 are now rows 10-12; delete them too.
$ws.Range("A10:A12").EntireRow.Delete()

# Renumber the "No" column for the remaining data rows.
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
